# Add new columns I (I0) and J (IF) to Sheet1, matching the header style
# used by the existing columns (B1:H1), and populate data rows 2-71.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header cells: copy the formatting from H1 (bold, bordered, centered) ---
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)
$ws.Application.CutCopyMode = $false

$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# --- Data cells: I2:I71 and J2:J71 ---
$I0 = @(9,7,8,1,8,8,7,9,7,6,1,9,9,8,7,5,8,7,5,7,7,5,8,7,7,7,8,7,7,8,7,7,6,7,6,7,7,9,10,6,9,5,6,5,8,6,9,4,5,7,7,10,8,5,9,6,4,5,6,5,6,1,7,5,7,7,8,7,9,3)
$IF = @(9,7,8,2,8,8,7,9,7,6,2,9,9,8,7,6,8,7,5,7,8,6,8,7,7,8,8,8,7,8,8,7,8,7,8,7,7,9,10,7,9,6,6,7,8,6,9,6,7,8,8,10,9,6,9,6,5,6,7,7,6,3,8,6,8,8,9,7,9,4)

for ($i = 0; $i -lt 70; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 9).Value = $I0[$i]
    $ws.Cells.Item($row, 10).Value = $IF[$i]
}
